$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 7618.0835
$ws.Range("I62").Value = 6849.8
$ws.Range("J62").Value = 11459.5
$ws.Range("K62").Value = 6849.8
$ws.Range("L62").Value = 11459.5
$ws.Range("M62").Value = -6225.8
$ws.Range("N62").Value = -12707.5

# Row 65
$ws.Range("H65").Value = 7618.0835
$ws.Range("I65").Value = 6849.8
$ws.Range("J65").Value = 11459.5
$ws.Range("K65").Value = 34249
$ws.Range("L65").Value = 57297.5
$ws.Range("M65").Value = -31129
$ws.Range("N65").Value = -63537.5

# Row 132
$ws.Range("H132").Value = 3492.342
$ws.Range("I132").Value = 2359.1177
$ws.Range("K132").Value = 7077.353099999999
$ws.Range("M132").Value = -4547.353099999999

# Row 137
$ws.Range("H137").Value = 1622.0385
$ws.Range("I137").Value = 1376.2778
$ws.Range("K137").Value = 4128.8334
$ws.Range("M137").Value = -1578.8334

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 5193.8184
$ws.Range("I74").Value = 3058
$ws.Range("K74").Value = 3058
$ws.Range("M74").Value = -2184

# Row 77
$ws.Range("H77").Value = 5193.8184
$ws.Range("I77").Value = 3058
$ws.Range("K77").Value = 15290
$ws.Range("M77").Value = -10922

# Row 110
$ws.Range("H110").Value = 1171.0714
$ws.Range("I110").Value = 1033
$ws.Range("J110").Value = 1999.5
$ws.Range("K110").Value = 1033
$ws.Range("L110").Value = 1999.5
$ws.Range("M110").Value = 1012
$ws.Range("N110").Value = -6089.5

# Row 132
$ws.Range("H132").Value = 1788.6136
$ws.Range("I132").Value = 1326.1666
$ws.Range("K132").Value = 3978.4998
$ws.Range("M132").Value = -1448.4998

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 745.069
$ws.Range("I7").Value = 853.05884
$ws.Range("K7").Value = 853.05884
$ws.Range("M7").Value = -740.05884

# Row 31
$ws.Range("H31").Value = 7229.204
$ws.Range("I31").Value = 3805.3
$ws.Range("J31").Value = 9590.518
$ws.Range("K31").Value = 3805.3
$ws.Range("L31").Value = 9590.518
$ws.Range("M31").Value = -3510.3
$ws.Range("N31").Value = -10180.518

# Row 34
$ws.Range("H34").Value = 7229.204
$ws.Range("I34").Value = 3805.3
$ws.Range("J34").Value = 9590.518
$ws.Range("K34").Value = 3805.3
$ws.Range("L34").Value = 9590.518
$ws.Range("M34").Value = -3603.3
$ws.Range("N34").Value = -9994.518

# Row 58
$ws.Range("H58").Value = 2746.95
$ws.Range("I58").Value = 2091.2666
$ws.Range("J58").Value = 4714
$ws.Range("K58").Value = 2091.2666
$ws.Range("L58").Value = 4714
$ws.Range("M58").Value = -1888.2666
$ws.Range("N58").Value = -5120

# Row 99
$ws.Range("H99").Value = 4856.125
$ws.Range("I99").Value = 4873.2
$ws.Range("J99").Value = 4827.6665
$ws.Range("K99").Value = 4873.2
$ws.Range("L99").Value = 4827.6665
$ws.Range("M99").Value = -3375.2
$ws.Range("N99").Value = -7823.6665

# Row 126
$ws.Range("H126").Value = 4856.125
$ws.Range("I126").Value = 4873.2
$ws.Range("J126").Value = 4827.6665
$ws.Range("K126").Value = 14619.6
$ws.Range("L126").Value = 14482.9995
$ws.Range("M126").Value = -12149.6
$ws.Range("N126").Value = -19422.9995

# Row 136
$ws.Range("H136").Value = 2746.95
$ws.Range("I136").Value = 2091.2666
$ws.Range("J136").Value = 4714
$ws.Range("K136").Value = 6273.7998
$ws.Range("L136").Value = 14142
$ws.Range("M136").Value = -3723.7998
$ws.Range("N136").Value = -19242

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 645
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 645
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1935
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -2159

# Row 7
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 17
$ws.Range("H17").Value = 1080.9
$ws.Range("J17").Value = 1150
$ws.Range("L17").Value = 3450
$ws.Range("N17").Value = -3788

# Row 34
$ws.Range("H34").Value = 9734.9
$ws.Range("J34").Value = 10749.777
$ws.Range("L34").Value = 32249.331
$ws.Range("N34").Value = -32417.331

# Row 39
$ws.Range("H39").Value = 9168.75
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 9764.286
$ws.Range("K39").Value = 15000
$ws.Range("L39").Value = 29292.858
$ws.Range("M39").Value = -14706
$ws.Range("N39").Value = -29880.858

# Row 55
$ws.Range("H55").Value = 9999.416999999999
$ws.Range("J55").Value = 9999.416999999999
$ws.Range("L55").Value = 29998.251
$ws.Range("N55").Value = -30352.251

# Row 63
$ws.Range("H63").Value = 15655.5
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 15655.5
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# Row 81
$ws.Range("H81").Value = 3000
$ws.Range("I81").Value = 3000
$ws.Range("K81").Value = 9000
$ws.Range("M81").Value = -7877

# Row 84
$ws.Range("H84").Value = 3000
$ws.Range("I84").Value = 3000
$ws.Range("K84").Value = 27000
$ws.Range("M84").Value = -21384

# Row 107
$ws.Range("H107").Value = 1483.6923
$ws.Range("J107").Value = 339.16666
$ws.Range("L107").Value = 1017.49998
$ws.Range("N107").Value = -4857.49998

# Row 122
$ws.Range("H122").Value = 1035.2
$ws.Range("I122").Value = 1498.5
$ws.Range("J122").Value = 726.3333
$ws.Range("K122").Value = 13486.5
$ws.Range("L122").Value = 6536.9997
$ws.Range("M122").Value = -11036.5
$ws.Range("N122").Value = -11436.9997

# Row 132
$ws.Range("H132").Value = 1487.0541
$ws.Range("I132").Value = 1514.8857
$ws.Range("K132").Value = 13633.9713
$ws.Range("M132").Value = -11103.9713

# Row 133
$ws.Range("H133").Value = 11830.77
$ws.Range("I133").Value = 4800
$ws.Range("J133").Value = 12416.667
$ws.Range("K133").Value = 14400
$ws.Range("L133").Value = 37250.001
$ws.Range("M133").Value = -9340
$ws.Range("N133").Value = -47370.001

# Row 134
$ws.Range("H134").Value = 4748.5
$ws.Range("I134").Value = 5464
$ws.Range("J134").Value = 4033
$ws.Range("K134").Value = 16392
$ws.Range("L134").Value = 12099
$ws.Range("M134").Value = -11322
$ws.Range("N134").Value = -22239

# Row 135
$ws.Range("H135").Value = 645
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 645
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 5805
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -10875

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 96
$ws.Range("H96").Value = 40250
$ws.Range("J96").Value = 40250
$ws.Range("L96").Value = 40250
$ws.Range("N96").Value = -45742

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Row 22
$ws.Range("H22").Value = 1970.375
$ws.Range("J22").Value = 1798.6
$ws.Range("L22").Value = 1798.6
$ws.Range("N22").Value = -2388.6

# Row 27
$ws.Range("H27").Value = 1970.375
$ws.Range("J27").Value = 1798.6
$ws.Range("L27").Value = 1798.6
$ws.Range("N27").Value = -2012.6

# Row 55
$ws.Range("H55").Value = 2528.4285
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 2866.5
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 2866.5
$ws.Range("N55").Value = -3212.5
$ws.Range("M55").Value = -327

# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 122
$ws.Range("H122").Value = 4261.3687
$ws.Range("I122").Value = 3574.5356
$ws.Range("K122").Value = 10723.6068
$ws.Range("M122").Value = -8273.606800000001

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2574.3142
$ws.Range("I122").Value = 2035.5483
$ws.Range("K122").Value = 6106.644899999999
$ws.Range("M122").Value = -3656.644899999999

# Row 132
$ws.Range("H132").Value = 2705.0557
$ws.Range("I132").Value = 2268.0571
$ws.Range("K132").Value = 6804.1713
$ws.Range("M132").Value = -4274.1713
